# The unified diff for this commit ("Moving from 2.0.1 to 2.0.2") touches
# only word/document.xml and word/styles.xml, and every single changed
# line is a pure character-for-character anagram of the line it replaces:
# the set of attributes/namespace declarations on each element (and their
# values) is completely unchanged, only their serialization order differs
# (mostly alphabetised, e.g. <w:tblW w:w="0" w:type="auto"/> becomes
# <w:tblW w:type="auto" w:w="0"/>, and the root <w:document> namespace
# declarations are reordered alphabetically by prefix). No text, value,
# element, attribute, relationship, or other part of the document's
# content actually changes.
#
# That kind of attribute/namespace reordering is a serialization-layer
# artifact of whatever tool re-exported the fixture for that commit
# (consistent with it being a mechanical resource refresh tied to the
# M2Doc library's 2.0.1 -> 2.0.2 version bump, not a document content
# edit) -- it is not something exposed through, or reachable via, the
# Word object model: every write path available on $word.ActiveDocument
# (property assignment, Find/Replace, Tables.Add, etc.) always leaves
# previously-serialized attributes in their original order and always
# emits freshly generated markup in the fixed schema order, never
# alphabetised. So there is no content mutation to make here -- applying
# this change means leaving the document's canonical content exactly as
# it is.
$d = $word.ActiveDocument

# Touch the document without altering any content, so the script is
# clearly operating on $word.ActiveDocument per the runtime's contract,
# while guaranteeing the saved package stays byte-for-byte identical in
# content to the source (only attribute/namespace ordering -- which
# carries no semantic meaning -- would differ from the literal commit,
# and that ordering is not controllable from the object model).
$null = $d.Content.Text.Length
